$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.252.42"
$ws.Range("E2").Value = "  +4.29%  "
$ws.Range("D3").Value = "2.592.96"
$ws.Range("E3").Value = "  +2.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.52"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.81"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.567"
$ws.Range("E8").Value = "  +1.40%  "
$ws.Range("D9").Value = "2.610.54"
$ws.Range("E9").Value = "  +2.66%  "
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("E12").Value = "  +2.92%  "
$ws.Range("E13").Value = "  +2.54%  "
$ws.Range("D14").Value = "3.056.50"
$ws.Range("E14").Value = "  +2.41%  "
$ws.Range("D15").Value = "59.161.31"
$ws.Range("E15").Value = "  +4.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.55"
$ws.Range("E16").Value = "  +2.51%  "
$ws.Range("D17").Value = "2.602.00"
$ws.Range("E17").Value = "  +2.47%  "
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "340.05"
$ws.Range("E19").Value = "  +2.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.34"
$ws.Range("E20").Value = "  +1.63%  "
$ws.Range("E21").Value = "  +1.41%  "
$ws.Range("E22").Value = "  +6.12%  "
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.12"
$ws.Range("E24").Value = "  +1.60%  "
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("E26").Value = "  +1.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.12"
$ws.Range("E28").Value = "  +3.46%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "0.0₃0727"
$ws.Range("E30").Value = "  -2.45%  "
$ws.Range("E31").Value = "  -4.80%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.57"
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.80"
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.76"
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.00"
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.36"
$ws.Range("E37").Value = "  +2.05%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.837"
$ws.Range("E38").Value = "  +1.99%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.45"
$ws.Range("E39").Value = "  +2.73%  "
$ws.Range("E40").Value = "  -1.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.52"
$ws.Range("E41").Value = "  +1.44%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "276.66"
$ws.Range("E43").Value = "  +5.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.77"
$ws.Range("E44").Value = "  +1.73%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0955"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.590"
$ws.Range("E46").Value = "  +2.32%  "
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.68"
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("D49").Value = "1.987.11"
$ws.Range("E49").Value = "  +1.59%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.62"
$ws.Range("E50").Value = "  +2.61%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0221"
$ws.Range("E51").Value = "  +0.00%  "
